$d = $word.ActiveDocument

# --- Change 3: remove the old _GoBack bookmark that sits next to the
#     page_total / page_total_master0 bookmarks, then re-add a fresh
#     _GoBack bookmark at the very start of the document (change 1).
#     Doing the delete first keeps the two still-existing bookmarks'
#     behaviour untouched, and the later Add mints a new id for the
#     relocated _GoBack bookmark.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$startRange = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $startRange)

# --- Change 2: split the run containing "3.1" into two runs, "3." and
#     "2", keeping identical run formatting on both pieces. Locate the
#     text first so this keeps working even if earlier edits shifted
#     character offsets.
$find = $d.Content.Find
$find.ClearFormatting()
[void]$find.Execute("3.1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $find.Parent.Start
$matchEnd = $find.Parent.End

# Replace "1" with "2" first (keeps a single run for now).
$tail = $d.Range($matchStart + 2, $matchEnd)
$tail.Text = "2"

# Force a genuine run boundary between "3." and "2" by nudging the
# second character's direct formatting away and then back again -
# Word (and this host) only materialises two separate <w:r> runs once
# the formatting has actually diverged at some point.
$tailAgain = $d.Range($matchStart + 2, $matchStart + 3)
$tailAgain.Font.Bold = $true
$tailAgain2 = $d.Range($matchStart + 2, $matchStart + 3)
$tailAgain2.Font.Bold = $false

Write-Output "done"
